$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "cfast pressure" values for K12:K101 (simulation re-run with longer time base)
$kValues = @(3219.29,6439.18,9659.66,12880.8,16102.4,19324.7,22547.6,25771.1,28995.2,32219.9,35445.2,38671.1,41897.6,45124.7,48352.4,51580.7,54809.6,58039.1,61269.1,64499.8,67731.1,70963,74195.5,77428.6,80662.3,83896.6,87131.5,90367,93603.1,96839.8,100077,103315,106553,109793,113032,116272,119513,122755,125997,129240,132483,135727,138971,142216,145462,148708,151955,155203,158451,161699,164949,168199,171449,174700,177952,181204,184457,187710,190964,194219,197474,200730,203987,207244,210501,213759,217018,220278,223538,226798,230060,233321,236584,239847,243111,246375,249640,252905,256171,259438,262705,265973,269241,272510,275780,279050,282321,285592,288864,292136)

# New "cfast temperature" values for G12:G101 (simulation re-run with longer time base)
$gValues = @(29.2764,38.5521,47.8271,57.1013,66.3748,75.6476,84.9197,94.191,103.462,112.731,122.001,131.269,140.537,149.804,159.07,168.335,177.6,186.864,196.128,205.39,214.652,223.913,233.174,242.433,251.692,260.95,270.208,279.465,288.721,297.976,307.231,316.484,325.737,334.99,344.241,353.492,362.743,371.992,381.241,390.489,399.736,408.983,418.228,427.473,436.718,445.961,455.204,464.446,473.688,482.929,492.169,501.408,510.646,519.884,529.121,538.358,547.593,556.828,566.062,575.296,584.528,593.76,602.991,612.222,621.452,630.681,639.909,649.137,658.363,667.59,676.815,686.04,695.264,704.487,713.709,722.931,732.152,741.372,750.592,759.811,769.029,778.246,787.463,796.679,805.894,815.109,824.322,833.535,842.748,851.959)

$kArr = New-Object 'object[,]' $kValues.Length,1
$gArr = New-Object 'object[,]' $gValues.Length,1
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $kArr[$i,0] = $kValues[$i]
    $gArr[$i,0] = $gValues[$i]
}

$ws.Range("K12:K101").Value = $kArr
$ws.Range("G12:G101").Value = $gArr

# Reflect the selection left after pasting in the updated pressure column
$ws.Range("K11:K101").Select()
